{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Find the \"Prueba 6\" Heading 2 paragraph that starts the block to remove.\nlet startIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.style === \"Heading 2\" && p.text.trim() === \"Prueba 6\") {\n    startIndex = i;\n    break;\n  }\n}\n\nif (startIndex === -1) {\n  throw new Error(\"Could not find the 'Prueba 6' heading paragraph\");\n}\n\n// Remove every paragraph from \"Prueba 6\" to the end of the body (this\n// includes the trailing empty paragraph right before the section break).\nfor (let i = paragraphs.items.length - 1; i >= startIndex; i--) {\n  paragraphs.items[i].delete();\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"Prueba 6\" Heading-2 paragraph (the start of the block to remove)\n# and the trailing empty paragraph right before the section break (the end of\n# the block to remove), then delete the whole range in one shot.\n\n$startPara = $null\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -eq \"Prueba 6\") {\n        $startPara = $p\n        break\n    }\n}\n\nif ($startPara -eq $null) {\n    throw \"Could not find the 'Prueba 6' heading paragraph\"\n}\n\n# The section's last paragraph is the empty one right before the sectPr;\n# that is simply the document's very last paragraph.\n$endPara = $d.Paragraphs.Last\n\n$range = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$range.Delete()\n"}
